$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")
$ws.Activate()
$ws.Range("B5").Value = "http://www.csm-testcenter.org/test?do=show&subdo=common&test=file_upload"
$ws.Range("B5").Select()
